$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row at position 882, shifting existing row 882 (and below) down by one.
$ws.Rows.Item(882).Insert()

# Populate the newly inserted row 882 with the new data point.
# Force column A to stay plain text (matches the sheet's other date cells,
# which are stored as literal strings, not real date serials), then restore
# the cell's format to match its plain, unstyled neighbours.
$ws.Cells.Item(882, 1).NumberFormat = "@"
$ws.Cells.Item(882, 1).Value = "2026/03/01"
$ws.Cells.Item(882, 1).Style = $ws.Cells.Item(881, 1).Style
$ws.Cells.Item(882, 2).Value = "日"
$ws.Cells.Item(882, 3).Value = 7
$ws.Cells.Item(882, 4).Value = 38
